$d = $word.ActiveDocument

$find = "V roku Súhvezdie Pegasus 2022: 8. – 17. október, 7. – 16. november,"
$replace = "V roku 2022 môžete pozorovať súhvezdie Súhvezdie Pegasus: 8. – 17. október, 7. – 16. november,"

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $replace, 2)
